# Update handback/handoff timestamps for the 7d100aed*.md row across sheets,
# as part of regenerating the Handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 7d100aed...md (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 04:46:44"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for row 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 04:46:40"
$wsZhCn.Range("K4").Value = "2016-08-31 04:46:58"

# de-de sheet: "Correspond Handback DateTime" for row 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-31 04:47:12"
